# Generate Report for Handoff
# A new handoff run replaces the old tracking GUID/hash pair
# (334cfc22-ec23-446b-86b3-a67debf0029e / 405231c618d8d6552b65d1ec475a33c0232764c8)
# with a new one, and records the new handoff timestamps, on all three worksheets.

$wb = $excel.ActiveWorkbook

$newGuid = "94f631a3-6c80-40bb-b9f5-46cbfbfaa00e"
$newHash = "cd78c68466b800bb972e08fd6d1601a0d6ee94b8"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Source markdown file name (column A, row 2) on every sheet.
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("A2").Value = "$newGuid.md"

# zh-cn handoff: new handoff file name + new handoff datetime.
$wsZhCn.Range("C2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-03-02 15:22:26"

# de-de handoff: new handoff file name + new handoff datetime.
$wsDeDe.Range("C2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-03-02 15:22:37"
